# Auto-generated script to apply scheduled market-data refresh to Goblin_Profits workbook.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H,I,J,K,L,M,N) on several rows
# across all eight crafting-job sheets, reflecting the latest market board data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2286.524
$ws.Range("J40").Value = 2478.889
$ws.Range("L40").Value = 2478.889
$ws.Range("N40").Value = -2828.889
$ws.Range("H74").Value = 7197.1665
$ws.Range("I74").Value = 5136.6
$ws.Range("K74").Value = 5136.6
$ws.Range("M74").Value = -4200.6
$ws.Range("H77").Value = 7197.1665
$ws.Range("I77").Value = 5136.6
$ws.Range("K77").Value = 25683
$ws.Range("M77").Value = -21003
$ws.Range("H88").Value = 11249.875
$ws.Range("I88").Value = 9999
$ws.Range("J88").Value = 11428.571
$ws.Range("K88").Value = 9999
$ws.Range("L88").Value = 11428.571
$ws.Range("M88").Value = -9593
$ws.Range("N88").Value = -12240.571
$ws.Range("H91").Value = 11249.875
$ws.Range("I91").Value = 9999
$ws.Range("J91").Value = 11428.571
$ws.Range("K91").Value = 9999
$ws.Range("L91").Value = 11428.571
$ws.Range("M91").Value = -8595
$ws.Range("N91").Value = -14236.571
$ws.Range("H116").Value = 5967.8
$ws.Range("I116").Value = 6922.25
$ws.Range("J116").Value = 2150
$ws.Range("K116").Value = 6922.25
$ws.Range("L116").Value = 2150
$ws.Range("M116").Value = -3480.25
$ws.Range("N116").Value = -9034
$ws.Range("H132").Value = 2945.842
$ws.Range("I132").Value = 2115.2942
$ws.Range("K132").Value = 6345.882599999999
$ws.Range("M132").Value = -3815.882599999999
$ws.Range("H137").Value = 9608.5
$ws.Range("I137").Value = 12679.75
$ws.Range("K137").Value = 38039.25
$ws.Range("M137").Value = -35489.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2672.889
$ws.Range("I32").Value = 2742.22
$ws.Range("K32").Value = 2742.22
$ws.Range("M32").Value = -2455.22
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H45").Value = 3963.8572
$ws.Range("J45").Value = 5486.75
$ws.Range("L45").Value = 5486.75
$ws.Range("N45").Value = -6240.75
$ws.Range("H102").Value = 2702.353
$ws.Range("I102").Value = 2702.353
$ws.Range("K102").Value = 2702.353
$ws.Range("M102").Value = -1080.353
$ws.Range("H132").Value = 3219.2354
$ws.Range("I132").Value = 3326.889
$ws.Range("K132").Value = 9980.667000000001
$ws.Range("M132").Value = -7450.667000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 67982
$ws.Range("J109").Value = 67982
$ws.Range("L109").Value = 67982
$ws.Range("N109").Value = -70756

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1915
$ws.Range("I132").Value = 1915
$ws.Range("K132").Value = 5745
$ws.Range("M132").Value = -3215
$ws.Range("H134").Value = 2932.9167
$ws.Range("J134").Value = 2814.4
$ws.Range("L134").Value = 8443.200000000001
$ws.Range("N134").Value = -13513.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 181.93333
$ws.Range("J12").Value = 205.27272
$ws.Range("L12").Value = 615.81816
$ws.Range("N12").Value = -961.81816
$ws.Range("H29").Value = 250625.25
$ws.Range("I29").Value = 500250
$ws.Range("K29").Value = 1500750
$ws.Range("M29").Value = -1500473
$ws.Range("H32").Value = 666806.7
$ws.Range("I32").Value = 500210
$ws.Range("J32").Value = 1000000
$ws.Range("K32").Value = 1500630
$ws.Range("L32").Value = 3000000
$ws.Range("M32").Value = -1500347
$ws.Range("N32").Value = -3000566
$ws.Range("H81").Value = 7089.6
$ws.Range("J81").Value = 8428.286
$ws.Range("L81").Value = 25284.858
$ws.Range("N81").Value = -27530.858
$ws.Range("H84").Value = 7089.6
$ws.Range("J84").Value = 8428.286
$ws.Range("L84").Value = 75854.57399999999
$ws.Range("N84").Value = -87086.57399999999
$ws.Range("H98").Value = 1839.6
$ws.Range("J98").Value = 1799.75
$ws.Range("L98").Value = 5399.25
$ws.Range("N98").Value = -8395.25
$ws.Range("H129").Value = 2947
$ws.Range("J129").Value = 3609.6667
$ws.Range("L129").Value = 10829.0001
$ws.Range("N129").Value = -20829.0001
$ws.Range("H131").Value = 3563995.2
$ws.Range("J131").Value = 4169643
$ws.Range("L131").Value = 12508929
$ws.Range("N131").Value = -12519009
$ws.Range("H140").Value = 69607.664
$ws.Range("I140").Value = 74293.92999999999
$ws.Range("K140").Value = 222881.79
$ws.Range("M140").Value = -217701.79

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4347.067
$ws.Range("I97").Value = 1162.5834
$ws.Range("K97").Value = 1162.5834
$ws.Range("M97").Value = -666.5834
$ws.Range("H126").Value = 3733.2222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 23991
$ws.Range("I43").Value = 22999
$ws.Range("J43").Value = 24983
$ws.Range("K43").Value = 22999
$ws.Range("L43").Value = 24983
$ws.Range("M43").Value = -22806
$ws.Range("N43").Value = -25369
$ws.Range("H61").Value = 5588.722
$ws.Range("J61").Value = 7508.909
$ws.Range("L61").Value = 7508.909
$ws.Range("N61").Value = -7912.909
$ws.Range("H93").Value = 3611.45
$ws.Range("I93").Value = 1873.1
$ws.Range("J93").Value = 5349.8
$ws.Range("K93").Value = 1873.1
$ws.Range("L93").Value = 5349.8
$ws.Range("M93").Value = -625.0999999999999
$ws.Range("N93").Value = -7845.8
$ws.Range("H100").Value = 3672.3635
$ws.Range("I100").Value = 3174.5
$ws.Range("K100").Value = 3174.5
$ws.Range("M100").Value = -2633.5
$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 40000
$ws.Range("L104").Value = 40000
$ws.Range("N104").Value = -46988
$ws.Range("H113").Value = 5588.722
$ws.Range("J113").Value = 7508.909
$ws.Range("L113").Value = 7508.909
$ws.Range("N113").Value = -11848.909

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7569.5
$ws.Range("I96").Value = 7569.5
$ws.Range("K96").Value = 7569.5
$ws.Range("M96").Value = -6196.5
$ws.Range("H132").Value = 2309.7273
$ws.Range("I132").Value = 2007.7
$ws.Range("K132").Value = 6023.1
$ws.Range("M132").Value = -3493.1
$ws.Range("H135").Value = 72999.664
$ws.Range("J135").Value = 72999.664
$ws.Range("L135").Value = 72999.664
$ws.Range("N135").Value = -83139.664

